# Commit: "update progress and action"
#
# 1. Rename the sheet "Activity Log" -> "Sheet1"
# 2. Give the header row (A1:F1) a bold font, a thin box border, and
#    center/top alignment
# 3. Move three messages from "Not started" to "In progress" and one to
#    "Done" in the Status column
# 4. Drop the Status column's old data-validation dropdown list

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "Sheet1"

# 2. Format the header row
$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous
$header.Borders.Weight = 2            # xlThin

# 3. Update the Status values
$ws.Range("F2").Value = "In progress"
$ws.Range("F3").Value = "In progress"
$ws.Range("F4").Value = "In progress"
$ws.Range("F5").Value = "Done"

# 4. Remove the old data-validation dropdown on the Status column
$ws.Range("F2:F1048576").Validation.Delete()
